# Fruta / hortaliza, semanal
# A new weekly price record (Higo, Mercado Mayorista Lo Valledor de Santiago)
# is inserted as row 31, pushing the previously existing rows 31-46 down to
# rows 32-47 (sheet grows from A1:T46 to A1:T47).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 31; this shifts rows 31..46 down
# to 32..47 and extends the used range to row 47.
$ws.Rows.Item(31).Insert()

# Populate the newly inserted row 31 with the new weekly record.
$ws.Range("A31").Value = 6
$ws.Range("B31").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C31").Value = "Metropolitana"
$ws.Range("D31").Value = 45029
$ws.Range("E31").Value = 13
$ws.Range("F31").Value = "Fruta"
$ws.Range("G31").Value = 100101
$ws.Range("H31").Value = "Berries"
$ws.Range("I31").Value = 100101006
$ws.Range("J31").Value = "Higo"
$ws.Range("K31").Value = "Sin especificar"
$ws.Range("L31").Value = "Primera"
$ws.Range("M31").Value = 120
$ws.Range("N31").Value = 20000
$ws.Range("O31").Value = 20000
$ws.Range("P31").Value = 20000
$ws.Range("Q31").Value = "$/bandeja 7 kilos"
$ws.Range("R31").Value = "Región Metropolitana"
$ws.Range("S31").Value = 2857
$ws.Range("T31").Value = 7
